$d = $word.ActiveDocument

# --- 1) New bullet: "Quand un utilisateur est supprimé..." right after
#        "Un forum d'entraide / de discussion" (was an empty paragraph) ---
$pForum = $d.Paragraphs.Item(25)
$pQuand = $d.Paragraphs.Item(26)
$pQuand.Range.InsertAfter("Quand un utilisateur est supprimé alors ses commentaires restent et le nom de l'utilisateur qui a posté le commentaire change. Par contre un commentaire peut avoir un ou plusieurs commentaire(s) en réponse. ")
$pQuand.Style = "Paragraphedeliste"
$pQuand.Range.ListFormat.ApplyListTemplateWithLevel($pForum.Range.ListFormat.ListTemplate, $true, 2, $false, $false)

# --- 2) Turn the "Comme je sais..." paragraph (text unchanged) into a
#        bulleted list paragraph as well ---
$pComme = $d.Paragraphs.Item(27)
$pComme.Style = "Paragraphedeliste"
$pComme.Range.ListFormat.ApplyListTemplateWithLevel($pQuand.Range.ListFormat.ListTemplate, $true, 2, $false, $false)

# --- 3) Remove the "Participer à des défis" bullet under "Utilisateur inscrit" ---
$rng = $d.Content
[void]$rng.Find.Execute("Participer à des défis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$rng.Expand(4)
$rng.Delete()

# --- 4) Fill the two empty trailing paragraphs at the end of the document:
#        one new bullet "Participer à des défis ou événements" and a
#        plain paragraph containing a single space ---
$count = $d.Paragraphs.Count
$pDefisEvt = $d.Paragraphs.Item($count - 1)
$pDefisEvt.Range.InsertAfter("Participer à des défis ou événements")
$pDefisEvt.Style = "Paragraphedeliste"
$pDefisEvt.Range.ListFormat.ApplyListTemplateWithLevel($pComme.Range.ListFormat.ListTemplate, $true, 2, $false, $false)

$pSpace = $d.Paragraphs.Item($count)
$pSpace.Range.InsertAfter(" ")

Write-Output "Done"
